# Fix scaling of 2024 logs: update column B (rows 1-48) with corrected values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(79, 100, 110, 127, 135, 180, 270, 288, 300, 322, 353, 376, 390, 440, 464, 497, 547, 563, 592, 596, 617, 655, 668, 691, 733, 766, 807, 846, 866, 944, 1013, 1036, 1080, 1136, 1185, 1214, 1222, 1187, 1323, 1344, 1379, 1388, 1453, 1505, 1534, 1552, 1581, 1625)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
